# Applies the edit described by the diff:
#  - Clears the RealLength(cm) values in S2:S75, S80, and S84:S101
#    (leaving already-blank cells such as S76:S79, S81:S83 untouched)
#  - Fills specific cells in S110:S133 with new computed RealLength(cm) values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear only the cells that previously held a numeric RealLength value.
$ws.Range("S2:S75").ClearContents()
$ws.Range("S80").ClearContents()
$ws.Range("S84:S101").ClearContents()

# Set new values for the specified rows in S110:S133
$newValues = @{
    110 = 114.7897603051682
    111 = 145.3017336543949
    113 = 106.6826406656907
    115 = 108.000663780192
    116 = 128.3856634920287
    121 = 114.4918344927808
    122 = 102.2261206162045
    127 = 119.1255119628378
    128 = 132.2738322968843
    129 = 156.2316338908785
    130 = 143.3123969655151
    131 = 98.78957849404372
    132 = 161.6013424872473
    133 = 141.7639407002064
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 19).Value = $newValues[$row]
}
